$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet table has a "virtual" header row (Column1..Column13) that
# previously lived on the table definition only (table ref started at row 2,
# reusing row 2's real header "Make"/"Model"/... as both the sheet header and
# the table header). Materialize that table header into an actual row 1,
# leaving row 2 ("Make", "Model", ...) untouched.
for ($i = 1; $i -le 13; $i++) {
    $ws.Cells.Item(1, $i).Value = "Column$i"
}

# Remove the stray "group separator" rows (Firewall(10), Router(37),
# Switch(9), Network Device(58)) which were merged A:M cells used purely as
# section headers. Unmerge then clear their contents so the rows collapse
# out of the saved sheetData entirely.
$groupRows = 3, 4, 15, 54
foreach ($r in $groupRows) {
    $rng = $ws.Range("A" + $r + ":M" + $r)
    $rng.UnMerge()
    $rng.ClearContents()
}

# Extend the Excel table so its header row is row 1 (matching the newly
# written Column1..Column13 header) instead of row 2.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:M78"))
